# "Costs not in M$ but in $"
# The O&M / CAPEX figures feeding the two log-log regressions on the
# "ReverseEngineer" sheet were being compared in $ against reference points
# that are actually in M$ (x1,000,000). Fix the two small tables (rows 9-11
# and rows 44-47) so the $ values used for the regression are expressed in
# plain dollars, matching the reference curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReverseEngineer")

# --- Table 1 (rows 2-16): G9:G11 / H9:H11 ---------------------------------
# G9:G11 used to just mirror H3:H5 (already $M); now convert to $.
$ws.Range("G9").Formula  = "=H3*1000000"
$ws.Range("G10").Formula = "=H4*1000000"
$ws.Range("G11").Formula = "=H5*1000000"

# H9:H11 (net of the 41% factor) recompute automatically off the new G values,
# but re-enter them explicitly so the formulas/values round-trip cleanly.
$ws.Range("H9").Formula  = "=G9*(1-0.41)"
$ws.Range("H10").Formula = "=G10*(1-0.41)"
$ws.Range("H11").Formula = "=G11*(1-0.41)"

# --- Table 2 (rows 37-52): G44:G47 ----------------------------------------
$ws.Range("G44").Formula = "=H38*1000000"
$ws.Range("G45").Formula = "=H39*1000000"
$ws.Range("G46").Formula = "=H40*1000000"
$ws.Range("G47").Formula = "=H41*1000000"
# H41 carried an explicit number format (MW-style) that must not leak onto
# G47 just because the new formula happens to reference it.
$ws.Range("G47").Style = "Normal"

# F14:F16 and F50:F52 (the LOG() columns feeding the charts) are unchanged
# formulas - they simply recompute from the corrected G/H values above.

# --- View state -------------------------------------------------------
# The workbook was left with the ReverseEngineer tab active/selected instead
# of Sheet1, scrolled down to the second table, with E50 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E50").Select()

# --- Chart 6 (second chart, CAPEX log-log plot) nudged slightly on canvas -
$co = $ws.ChartObjects(2)
$co.Top = $co.Top - 10
$co.Left = $co.Left + 6.4285826771653545
